$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format so numeric-looking strings (e.g. "6.81") are not
# auto-converted to numbers by Excel, then restore default "Normal" style so
# no visible formatting change is introduced.
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "E8", "D9", "E9", "E10", "D11", "E11", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "E31", "D32", "E32", "D34", "E34", "D35", "E35", "D36", "E36", "E37", "D38", "E38", "E39", "E40", "E41", "E42", "D43", "E43", "E44", "E45", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.032.01'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '2.304.09'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '301.06'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '98.38'
$ws.Range("E6").Value = '  -1.16%  '
$ws.Range("E7").Value = '  +2.17%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").Value = '6.81'
$ws.Range("E14").Value = '  -1.77%  '
$ws.Range("D15").Value = '2.662.51'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '2.295.87'
$ws.Range("E16").Value = '  -2.80%  '
$ws.Range("E17").Value = '  -1.91%  '
$ws.Range("D18").Value = '42.990.85'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '6.13'
$ws.Range("E21").Value = '  -1.43%  '
$ws.Range("D22").Value = '68.36'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '242.19'
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").Value = '25.20'
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("D29").Value = '166.95'
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("D30").Value = '2.04'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").Value = '33.32'
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("D34").Value = '4.80'
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("D35").Value = '5.03'
$ws.Range("E35").Value = '  -2.90%  '
$ws.Range("D36").Value = '17.76'
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = '0.0689'
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("E39").Value = '  -1.57%  '
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").Value = '1.998.87'
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("E46").Value = '  +1.18%  '
$ws.Range("D47").Value = '17.44'
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").Value = '2.81'
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").Value = '53.66'
$ws.Range("E49").Value = '  -2.50%  '
$ws.Range("D50").Value = '2.528.17'
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").Value = '72.86'
$ws.Range("E51").Value = '  -4.55%  '

foreach ($addr in $cells) {
    $ws.Range($addr).Style = "Normal"
}
